$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Posição"

$ws.Range("M22").Select()
